$wb = $excel.ActiveWorkbook

# --- INDICATOR sheet: rename indicator IDs (labels unchanged) ---
$ws = $wb.Worksheets.Item("INDICATOR")
$ws.Range("A2").Value = "NMGDP"
$ws.Range("A4").Value = "NMWGT"

# Column B width tweak (20.28515625 -> 21)
$ws.Columns.Item(2).ColumnWidth = 20.16

# Selection / active cell moves from F10 to A6, and this sheet is no longer the tab shown on open
$ws.Range("A6").Select() | Out-Null

# --- INDUSTRY_TYPE sheet becomes the active/selected tab ---
$ws4 = $wb.Worksheets.Item("INDUSTRY_TYPE")
$ws4.Activate() | Out-Null
$ws4.Range("B2").Select() | Out-Null
